$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Parent initial measurements added - update stress-level answers for rows 2-7
$ws.Range("B2").Value = "Moderately stressful"
$ws.Range("B3").Value = "Very stressful"
$ws.Range("B4").Value = "A little stressful "
$ws.Range("B5").Value = "Not stressful"
$ws.Range("B6").Value = "Moderately stressful"
$ws.Range("B7").Value = "Moderately stressful"
